# Daily attendance processing - 2026-01-19 11:38:08
# Swap the order of "System" and the email address in column G
# (the "Updated By" / modifier column) from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" on every row where it occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
